$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9765523076057434
$ws.Range("B1").Value = 2.722717523574829
$ws.Range("C1").Value = 4.720315933227539
$ws.Range("D1").Value = 1.194574475288391
$ws.Range("E1").Value = 1.294441223144531
